# Add season-record columns (Wins / Losses / Ties) to the player table.
# The previous scraper only pulled team statistics, not the season record,
# so this fills in the 2010 Houston Astros record (76-86-0) for every
# player row, matching the header/style conventions already used by the
# sheet (columns A..AC).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, one column to the right of the existing "Unnamed: 28"
# column (AC). Clone the formatting (bold font + thin border) from the
# neighboring header cell so the new headers look consistent with the rest
# of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row gets the same team season record: 76 wins, 86 losses,
# 0 ties.
for ($row = 2; $row -le 50; $row++) {
    $ws.Cells.Item($row, 30).Value = 76
    $ws.Cells.Item($row, 31).Value = 86
    $ws.Cells.Item($row, 32).Value = 0
}
